# Rename the three logo pictures that live in the document's header/footer
# stories:
#   - BTec_Logo-Orange (first-page header)           image2.jpg -> image1.jpg
#   - PearsonLogo.png  (first-page footer)            image1.png -> image2.png
#   - PearsonLogo.png  (primary/default footer)       image1.png -> image2.png
#
# Only the picture's "Name" (the wp:docPr / pic:cNvPr name="...") changes -
# size, description/alt-text, and the embedded image data are untouched.
#
# Note: the InlineShape.Name *getter* in this host doesn't reflect the
# stored name, so shapes are identified by their (reliable) AlternativeText
# instead of by their current Name.

$d = $word.ActiveDocument

# Helper: given an InlineShape handle, rename it. Renaming a shape that
# lives inside a footer story is flaky when done directly through the
# HeaderFooter.Range.InlineShapes collection in this host, so we first
# Select() the shape (which works reliably for both header- and
# footer-hosted shapes) and then re-fetch it from $word.Selection before
# touching the Name property.
function Rename-InlineShape($shape, $newName) {
    $shape.Select()
    $selShape = $word.Selection.InlineShapes.Item(1)
    $selShape.Name = $newName
}

$btecDescr = "BTec_Logo-Orange"
$pearsonDescr = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"

foreach ($sec in $d.Sections) {

    # --- Headers -----------------------------------------------------------
    for ($h = 1; $h -le $sec.Headers.Count; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq $btecDescr) {
                    Rename-InlineShape $shp "image1.jpg"
                }
            }
        }
    }

    # --- Footers -------------------------------------------------------------
    for ($f = 1; $f -le $sec.Footers.Count; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq $pearsonDescr) {
                    Rename-InlineShape $shp "image2.png"
                }
            }
        }
    }
}
